# Notes.xlsx: fix FirstName/LastName header order, rename Note_G3EI12x ->
# Note_G3EI11x columns, renumber the CNE (student id) column, move the
# selection, and re-save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap the FirstName/LastName labels (the data columns were
# already FirstName-then-LastName; the header text had them backwards) and
# rename the last three "Note_G3EI12x" headers to "Note_G3EI11x".
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"
$ws.Range("J1").Value = "Note_G3EI117"
$ws.Range("K1").Value = "Note_G3EI118"
$ws.Range("L1").Value = "Note_G3EI119"

# --- CNE (student id) column: give each student row its own sequential id
# instead of repeating the first student's id down the whole column.
$ws.Range("A3").Value = 19000042
$ws.Range("A4").Value = 19000043
$ws.Range("A5").Value = 19000044
$ws.Range("A6").Value = 19000045
$ws.Range("A7").Value = 19000046
$ws.Range("A8").Value = 19000047
$ws.Range("A9").Value = 19000048
$ws.Range("A10").Value = 19000049
$ws.Range("A11").Value = 19000050

# --- Move the active selection.
$ws.Range("H7").Select() | Out-Null
